$wb = $excel.ActiveWorkbook

# --- Remove the two data sheets (d-dataSheet, d-dev) ---
$wb.Worksheets("d-dataSheet").Delete()
$wb.Worksheets("d-dev").Delete()

# --- Add the new TestDataReader rows to the c-demo_ui config sheet ---
$ws1 = $wb.Worksheets("c-demo_ui")
$ws1.Range("A20").Value = "testdata.filename"
$ws1.Range("A21").Value = "testdata.sheetname"
$ws1.Range("B20").Value = "testdata.xlsx"

# --- Make c-demo_ui the active/selected sheet & tab (was t-soaptest before) ---
$ws1.Activate()
$ws1.Range("B21").Select()
